# Backlog.xlsx update — "Wykonanie zadan oraz aktualizacja dokumentacji"
#
# 1. Two more backlog tasks were finished by Przemek in sprint 3 (new table
#    rows 35-36), which grows the "Tabela3" listobject.
# 2. The per-person "3 sprint" totals (col I, rows 31-34) now sum over a
#    wider/unshared range and pick up the two new rows.
# 3. The old "calosc" (grand total) helper column that lived in R:S got
#    copied over to new columns V:W, and R:S was repurposed as the "3 sprint"
#    breakdown header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Two finished tasks for Przemek, sprint 3 ---------------------------
# Row 35
$ws.Range("A34").Copy()
$ws.Range("A35").PasteSpecial(-4122)          # xlPasteFormats: reuse the date style
$ws.Range("A35").Value2 = 45452
$ws.Range("B35").Value2 = "Przemek"
$ws.Range("C35").Value2 = "main: dodać funkcjonalność przycisku do wychodzenia z programu"
$ws.Range("D35").Value2 = 3
$ws.Range("E35").Value2 = 60
$ws.Range("F35").Value2 = "Ukończono"

# Row 36
$ws.Range("A34").Copy()
$ws.Range("A36").PasteSpecial(-4122)
$ws.Range("A36").Value2 = 45452
$ws.Range("B36").Value2 = "Przemek"
$ws.Range("C36").Value2 = "trening: dodanie graficznego interfejsu"
$ws.Range("D36").Value2 = 3
$ws.Range("E36").Value2 = 60
$ws.Range("F36").Value2 = "Ukończono"

# Grow the backlog table so the two new rows become part of it / the filter.
$lo = $ws.ListObjects.Item("Tabela3")
$lo.Resize($ws.Range("A2:F36"))

# --- 2. Re-point the sprint-3 per-person SUMIFs to cover the new rows ------
$ws.Range("I31").Formula = "=SUMIF(B6:B34,H3,E6:E107)/60"
$ws.Range("I32").Formula = "=SUMIF(B7:B35,H4,E7:E107)/60"
$ws.Range("I33").Formula = "=SUMIF(B8:B36,H5,E8:E107)/60"
$ws.Range("I34").Formula = "=SUMIF(B9:B37,H6,E9:E107)/60"

# --- 3. Move the "calosc" summary from R:S to V:W, reuse R:S for "3 sprint" -
$ws.Range("R1").Copy()
$ws.Range("V1").PasteSpecial(-4122)
$ws.Range("V1").Value2 = "całość"
$ws.Range("S1").Copy()
$ws.Range("W1").PasteSpecial(-4122)
$ws.Range("W1").Value2 = "h pracy"

for ($r = 3; $r -le 6; $r++) {
  $ws.Range("R$r").Copy()
  $ws.Range("V$r").PasteSpecial(-4122)
  $ws.Range("V$r").Value2 = $ws.Range("R$r").Value2

  $ws.Range("S$r").Copy()
  $ws.Range("W$r").PasteSpecial(-4122)
}

$ws.Range("R1").Value2 = "3 sprint"

# --- cosmetic: leave the selection near the newly-added rows ---------------
$ws.Range("I35").Select()
